# "Resultados tempos finais concluido"
# Fill in the previously-empty column I (rows 2-14) on Sheet1 with the
# measured timings, then leave Sheet1's selection on I13 and make Sheet2
# the active/selected tab (matching the tabSelected move + new selection
# recorded in the diff).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Column I timings (Sheet1, rows 2-14).
$ws1.Range("I2").Value  = "2.06 ms"
$ws1.Range("I3").Value  = " 0.10 ms"
$ws1.Range("I4").Value  = "59.42 ms"
$ws1.Range("I5").Value  = "0.21 ms"
$ws1.Range("I6").Value  = "0.75 ms"
$ws1.Range("I7").Value  = "0.09 ms"
$ws1.Range("I8").Value  = "0.07 ms"
$ws1.Range("I9").Value  = "0.14 ms"
$ws1.Range("I10").Value = "1.42 ms"
$ws1.Range("I11").Value = " 0.19 ms"
$ws1.Range("I12").Value = "88.71 ms"
$ws1.Range("I13").Value = "0.24 ms"
$ws1.Range("I14").Value = "0.02 ms"

# Sheet1's own selection moves to I13 (no longer the active/selected tab).
$ws1.Range("I13").Select()

# Sheet2 becomes the active/selected tab; its own selection (C13) is
# untouched.
$ws2.Activate()
